# Atualizacao rapida de agenda as  8:58:37,26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep only the first agenda entry (Pedro / 2590 / Sucata THL / ...),
# clearing out all the other previously-listed jobs (rows 3-8).
$ws.Range("A2").Value = "Pedro"
$ws.Range("B2").Value = "'2590"
$ws.Range("C2").Value = "Sucata THL"
$ws.Range("D2").Value = "Zona estava aberta, técnico vai verificar o que pode ser."
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = "Em andamento"
$ws.Range("H2").Value = "Maxvel: 38 / Forte: 17"
$ws.Range("I2").Value = $null

# Rows 3-8 get fully cleared (contents only, styles stay as-is)
$ws.Range("A3:I8").ClearContents()
$ws.Range("A3:I8").EntireRow.AutoFit()

# Update the view: scroll / selection moved
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("I5:I7").Select()
